$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: id=1, date=43363, type=EXPENSE, description="first salary", amount=25000
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 43363
$ws.Range("C2").Value = "EXPENSE"
$ws.Range("D2").Value = "first salary"
$ws.Range("E2").Value = 25000

# Row 3: id=2, date=43363, type=INCOME, description="eating lunch in the morning", amount=200
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 43363
$ws.Range("C3").Value = "INCOME"
$ws.Range("D3").Value = "eating lunch in the morning"
$ws.Range("E3").Value = 200

# Row 4: id=3, date=43363, type=EXPENSE, description="eating lunch", amount=500
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 43363
$ws.Range("C4").Value = "EXPENSE"
$ws.Range("D4").Value = "eating lunch"
$ws.Range("E4").Value = 500

# Row 5: id=4, date=43363, type=INCOME, description="salary", amount=20000
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 43363
$ws.Range("C5").Value = "INCOME"
$ws.Range("D5").Value = "salary"
$ws.Range("E5").Value = 20000

# Remove old row 6 (transaction deleted)
$ws.Range("A6:E6").Delete()
